# Commit: "Removed language defaults from template"
#
# The template previously hard-coded the survey/choices header labels as
# "label::language" / "hint::language". Those per-language defaults are
# dropped in favour of the plain "label" / "hint" headers.

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# survey sheet: C1 "label::language" -> "label", D1 "hint::language" -> "hint"
$survey.Range("C1").Value = "label"
$survey.Range("D1").Value = "hint"

# choices sheet: C1 "label::language" -> "label"
$choices.Range("C1").Value = "label"

# Reset the saved view/selection state on both sheets (as in the target
# workbook) without leaving "choices" as the active tab.
$choices.Range("A2").Select()
$survey.Activate()
$survey.Range("A2").Select()
